$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the missing Mass values for Ruthenium, Cymene and PTA
$ws.Range("D4").Value = 102.911626
$ws.Range("D8").Value = 135.116827
$ws.Range("D9").Value = 158.08416

# Reflect the last selected cell left by the author's editing session
$ws.Range("H20").Select()
